# Add a new "ListBullet" paragraph listing the responsible instructor
# right after the "Docente(s) Responsável(eis)" heading paragraph.

$d = $word.ActiveDocument

# Find the heading paragraph "Docente(s) Responsável(eis) " and obtain
# the Paragraph object that contains it.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Docente(s) Responsável(eis) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $headingPara = $searchRange.Paragraphs(1)

    # Insert a brand-new empty paragraph right after the heading.
    $headingPara.Range.InsertParagraphAfter()

    # The newly created paragraph is the one following the heading.
    $newPara = $headingPara.Next()
    $newPara.Range.Text = "6712818 - Mauricio Lamano Ferreira"
    $newPara.Style = "ListBullet"
}
